$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Row 2
$ws.Range("A2").Value = '2025-10-30 18:27:16'
$ws.Range("B2").Value = '大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5423720'
$ws.Range("G2").Value = 385
$ws.Range("H2").Value = '🔥AI,Ai ◆効率化'

# Row 3
$ws.Range("A3").Value = '2025-10-30 18:27:16'
$ws.Range("B3").Value = '【急募】映像解析AIによる自動検出・モザイク処理スクリプト開発'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5424032'
$ws.Range("G3").Value = 368
$ws.Range("H3").Value = '🔥AI,Ai ◆開発'

# Row 4
$ws.Range("A4").Value = '2025-10-30 18:27:16'
$ws.Range("B4").Value = 'Excel・Accessベースの改修や追加、Pythonスクレイピングやデータ整形等の開発員募集'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5405426'
$ws.Range("G4").Value = 305
$ws.Range("H4").Value = '🔥Python ◆開発,スクレイピング'

# Row 5
$ws.Range("A5").Value = '2025-10-30 18:27:16'
$ws.Range("B5").Value = '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5217096'
$ws.Range("G5").Value = 243
$ws.Range("H5").Value = '🔥API ◆ツール'

# Row 6
$ws.Range("A6").Value = '2025-10-30 18:27:16'
$ws.Range("B6").Value = '大手クレジットカード企業向け、Google Cloudを利用したアジャイル開発共通基盤案件'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5418320'
$ws.Range("G6").Value = 75
$ws.Range("H6").Value = '◆開発'

# Row 7
$ws.Range("A7").Value = '2025-10-30 18:27:16'
$ws.Range("B7").Value = '大手クレジットカード企業向け、Google Cloudを利用したアジャイル開発共通基盤案件_ワーカー'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5418318'
$ws.Range("G7").Value = 75
$ws.Range("H7").Value = '◆開発'

# Row 8
$ws.Range("A8").Value = '2025-10-30 18:27:16'
$ws.Range("B8").Value = 'wordpressレンダリングを妨げるリソースの除外'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5016989'
$ws.Range("G8").Value = 33
$ws.Range("H8").Value = '○WordPress'

# Row 9
$ws.Range("A9").Value = '2025-10-30 18:27:16'
$ws.Range("B9").Value = 'PaddlePaddle/PaddleOCR 文字列OCR 学習(検出+認識)'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '1,000,000 円 ~ 3,000,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5423522'
$ws.Range("G9").Value = 25

# Row 10
$ws.Range("A10").Value = '2025-10-30 18:27:16'
$ws.Range("B10").Value = '【フルリモート】SESエンジニア募集|スキルに応じて30〜40万円/月|複数案件あり・継続前提'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5417644'
$ws.Range("G10").Value = 25

# Row 11
$ws.Range("A11").Value = '2025-10-30 18:27:16'
$ws.Range("B11").Value = '評価基板設計・製造について'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5423728'
$ws.Range("G11").Value = 18

# Row 12
$ws.Range("A12").Value = '2025-10-30 18:27:16'
$ws.Range("B12").Value = '【急募】非接触センサー×Bluetoothデバイスのアイデア壁打ち相談'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5423605'
$ws.Range("G12").Value = 10

# Row 13
$ws.Range("A13").Value = '2025-10-30 18:27:16'
$ws.Range("B13").Value = 'Google Workspace × さくらサーバー DNSメール設定代行依頼'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5423476'
$ws.Range("G13").Value = 10

# Rebuild hyperlinks fresh (targets shifted along with the rows)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), $ws.Range("F2").Value())
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), $ws.Range("F3").Value())
$ws.Range("F3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F4"), $ws.Range("F4").Value())
$ws.Range("F4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F5"), $ws.Range("F5").Value())
$ws.Range("F5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F6"), $ws.Range("F6").Value())
$ws.Range("F6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F7"), $ws.Range("F7").Value())
$ws.Range("F7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F8"), $ws.Range("F8").Value())
$ws.Range("F8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F9"), $ws.Range("F9").Value())
$ws.Range("F9").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F10"), $ws.Range("F10").Value())
$ws.Range("F10").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F11"), $ws.Range("F11").Value())
$ws.Range("F11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F12"), $ws.Range("F12").Value())
$ws.Range("F12").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F13"), $ws.Range("F13").Value())
$ws.Range("F13").Style = "Hyperlink"

Write-Host "done"